$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing price values (K3:K6)
$ws.Range("K3").Value = 200
$ws.Range("K4").Value = 100
$ws.Range("K5").Value = 135
$ws.Range("K6").Value = 306

# Add a new row 7, duplicating the data pattern of row 2 (TD / Yevhen Khreptun)
# but with a new reservation id, new check-in/check-out dates, and a new price.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "TD"
$ws.Range("C7").Value = 203
$ws.Range("D7").Value = "td@test.ca"
$ws.Range("E7").Value = "Yevhen"
$ws.Range("F7").Value = "Khreptun"
$ws.Range("G7").Value = "yevhen@test.ca"
$ws.Range("H7").Value = 45211
$ws.Range("I7").Value = "24/12/2023"
$ws.Range("J7").Value = "30/12/2023"
$ws.Range("K7").Value = 250
$ws.Range("L7").Value = "confirmed"

# Match the date formatting used by the rest of the table (reuses the
# existing date-format style rather than minting a new one).
$ws.Range("H7").NumberFormat = "m/d/yy"
$ws.Range("I7").NumberFormat = "m/d/yy"
$ws.Range("J7").NumberFormat = "m/d/yy"

# Apply hyperlinks for the owner email (D7) and user email (G7) cells, matching
# the pattern used by the other rows.
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:td@test.ca")
$ws.Hyperlinks.Add($ws.Range("G7"), "mailto:yevhen@test.ca")

# Hyperlinks.Add mints its own style variant; reapply the shared hyperlink
# style so the new cells line up with the rest of the table.
$ws.Range("D7").Style = $ws.Range("D2").Style
$ws.Range("G7").Style = $ws.Range("G2").Style

# Update the selected cell to mirror the new active selection shown in the diff.
$ws.Range("F7").Select()
